$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the weight/volume (columns C/D) values for rows 2-26 back to their
# pre-bugfix-"fix" numbers (whole/round numbers instead of the fractional
# "corrected" ones), per the commit being reverted.
$values = @{
    2  = @(3, 1)
    3  = @(5, 1)
    4  = @(5, 1)
    5  = @(5, 1)
    6  = @(5, 1)
    7  = @(5, 1)
    8  = @(1, 1)
    9  = @(5, 1)
    10 = @(5, 1)
    11 = @(3, 1)
    12 = @(5, 1)
    13 = @(5, 1)
    14 = @(5, 1)
    15 = @(5, 1)
    16 = @(5, 1)
    17 = @(5, 1)
    18 = @(5, 1)
    19 = @(10, 2)
    20 = @(10, 2)
    21 = @(10, 2)
    22 = @(10, 2)
    23 = @(10, 2)
    24 = @(20, 3)
    25 = @(20, 3)
    26 = @(20, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("C$row").Value2 = $pair[0]
    $ws.Range("D$row").Value2 = $pair[1]
}

# Remove the stray K20 cell (and its now-unused shared string " ") that was
# added by the reverted commit.
$ws.Range("K20").ClearContents()

# Restore the prior selection.
[void]$ws.Range("G13").Select()
